# Weekly data refresh: a new price record (week of 2023-10-26) is inserted
# at the top of the data table (row 61), pushing all subsequent rows down
# by one. This mirrors how the source dataset is normally updated: the most
# recent week's observation is prepended above the existing history.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 61; everything from row 61 downward shifts
# down to row 62 onward (existing formatting/style on column D is carried
# along automatically by Excel's insert behavior).
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A61").Value = 7
$ws.Range("B61").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C61").Value = "Ñuble"
$ws.Range("D61").Value = 45225
$ws.Range("E61").Value = 16
$ws.Range("F61").Value = 100112001
$ws.Range("G61").Value = "Berenjena"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 30
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = 10000
$ws.Range("N61").Value = "$/caja 60 unidades"
$ws.Range("O61").Value = "Región de Arica y Parinacota"
$ws.Range("P61").Value = 167
$ws.Range("Q61").Value = 60
$ws.Range("R61").Value = "Hortaliza"
